$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial value for every data row (2-74).
# It currently stores serial 45172 (2023-09-03) and needs to become
# serial 45175 (2023-09-06) for all of these rows.
for ($row = 2; $row -le 74; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45172) {
        $cell.Value = 45175
    }
}
